# Update countries & provincias Spain
# Applies the daily refresh of the "Pais" (countries) COVID table:
#  - Re-sorted rows for three country pairs (text + stats swap positions)
#  - Updated totals for a handful of rows (including the two re-sorted pairs)
#  - Updated "last refreshed" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Country name swaps caused by re-sorting the table by "Casos totales"
# ---------------------------------------------------------------------------

# Mayotte / Tailandia swap rank (rows 140-141)
$ws.Range("A140").Value = "Tailandia"
$ws.Range("A141").Value = "Mayotte"

# Belice jumps ahead of Principado de Andorra / Nueva Zelanda (rows 156-158)
$ws.Range("A156").Value = "Belice"
$ws.Range("A157").Value = "Principado de Andorra"
$ws.Range("A158").Value = "Nueva Zelanda"

# Timor Oriental / Santa Lucia swap rank (rows 207-208)
$ws.Range("A207").Value = "Timor Oriental"
$ws.Range("A208").Value = "Santa Lucia"

# ---------------------------------------------------------------------------
# 2) Updated daily statistics (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes)
# ---------------------------------------------------------------------------

# Row 60 (Uzbekistan)
$ws.Range("B60").Value = 55593
$ws.Range("C60").Value = 273
$ws.Range("D60").Value = 51965
$ws.Range("E60").Value = 3168
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 460

# Row 66 (Kirguistan)
$ws.Range("B66").Value = 46355
$ws.Range("C66").Value = 104
$ws.Range("D66").Value = 42613
$ws.Range("E66").Value = 2678
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = 1064

# Row 140 (now Tailandia)
$ws.Range("B140").Value = 3545
$ws.Range("C140").Value = 22
$ws.Range("D140").Value = 3369
$ws.Range("E140").Value = 117
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 59

# Row 141 (now Mayotte)
$ws.Range("B141").Value = 3541
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 2964
$ws.Range("E141").Value = 537
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 40

# Row 156 (now Belice)
$ws.Range("B156").Value = 1854
$ws.Range("C156").Value = 29
$ws.Range("D156").Value = 1196
$ws.Range("E156").Value = 634
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 24

# Row 157 (now Principado de Andorra)
$ws.Range("B157").Value = 1836
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = 1263
$ws.Range("E157").Value = 520
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 53

# Row 158 (now Nueva Zelanda)
$ws.Range("B158").Value = 1833
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 1753
$ws.Range("E158").Value = 55
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 25

# ---------------------------------------------------------------------------
# 3) Refresh timestamp text
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Septiembre de 2020 a las 07:38"
